{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Change 1 & 2: collapse paragraphs that were split into multiple runs\n// (because of inline grammar-check markers) back into a single plain run,\n// by replacing each paragraph's content with its own full text. ---\nconst mergeTexts = [\n  \"Vorbau (abgesehen von Sensorik) abfeilen wo die schrauben drin sitzen, sodass es Rechteckig ist (keine Verschwendung\",\n  \"Vorteil: Effizient, Platzsparend, gleicher Wendekreis\",\n];\n\nlet lastParagraph = null;\nfor (const para of paragraphs.items) {\n  if (mergeTexts.includes(para.text)) {\n    para.insertText(para.text, Word.InsertLocation.replace);\n  }\n  lastParagraph = para;\n}\n\n// --- Change 3: append a new \"Einkaufsrat\" section at the end of the\n// document (after the final \"Nachteil: Schlechter Wendekreis\" paragraph,\n// before the section break). ---\nconst anchor = lastParagraph;\n\nconst p4 = anchor.insertParagraph(\"Nachteil: Nicht so sch\u00f6n\", Word.InsertLocation.after);\nconst p3 = anchor.insertParagraph(\"Vorteil: Sehr einfach\", Word.InsertLocation.after);\nconst p2 = anchor.insertParagraph(\n  \"Einkaufsrad unter den Vorbau, dadurch h\u00e4lft das Gewicht, problem gel\u00f6st.\",\n  Word.InsertLocation.after\n);\nconst p1 = anchor.insertParagraph(\"Einkaufsrat\", Word.InsertLocation.after);\np1.styleBuiltIn = Word.Style.heading3;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 & 2: collapse paragraphs that were split into multiple runs\n# (because of inline grammar-check markers) back into a single plain run,\n# by replacing each paragraph's own content with its own full text. ---\n$mergeTexts = @(\n    \"Vorbau (abgesehen von Sensorik) abfeilen wo die schrauben drin sitzen, sodass es Rechteckig ist (keine Verschwendung\",\n    \"Vorteil: Effizient, Platzsparend, gleicher Wendekreis\"\n)\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $r = $p.Range\n    $r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark (wdCharacter = 1)\n    $fullText = $r.Text\n    if ($mergeTexts -contains $fullText) {\n        $r.Delete()\n        $r.InsertAfter($fullText)\n    }\n}\n\n# --- Change 3: append a new \"Einkaufsrat\" section at the end of the\n# document (after the final \"Nachteil: Schlechter Wendekreis\" paragraph,\n# before the section break). ---\n$anchor = $d.Paragraphs($d.Paragraphs.Count)\n\n$r = $anchor.Range\n$r.Collapse(0) | Out-Null   # wdCollapseEnd\n$r.InsertParagraphAfter()\n$d.Range($anchor.Range.End, $anchor.Range.End).Paragraphs(1).Range.Text = \"Nachteil: Nicht so sch\u00f6n\"\n\n$r = $anchor.Range\n$r.Collapse(0) | Out-Null\n$r.InsertParagraphAfter()\n$d.Range($anchor.Range.End, $anchor.Range.End).Paragraphs(1).Range.Text = \"Vorteil: Sehr einfach\"\n\n$r = $anchor.Range\n$r.Collapse(0) | Out-Null\n$r.InsertParagraphAfter()\n$d.Range($anchor.Range.End, $anchor.Range.End).Paragraphs(1).Range.Text = \"Einkaufsrad unter den Vorbau, dadurch h\u00e4lft das Gewicht, problem gel\u00f6st.\"\n\n$r = $anchor.Range\n$r.Collapse(0) | Out-Null\n$r.InsertParagraphAfter()\n$heading = $d.Range($anchor.Range.End, $anchor.Range.End).Paragraphs(1)\n$heading.Range.Text = \"Einkaufsrat\"\n$heading.Style = $d.Styles(\"Heading 3\")\n"}
